# Auto-generated edit script: apply scheduled-runner price refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 607.4545000000001
$ws.Range("I6").Value = 607.4545000000001
$ws.Range("K6").Value = 1822.3635
$ws.Range("M6").Value = -1710.3635

$ws.Range("H8").Value = 189.66667
$ws.Range("I8").Value = 189.66667
$ws.Range("K8").Value = 569.00001
$ws.Range("M8").Value = -430.00001

$ws.Range("H31").Value = 103.75
$ws.Range("I31").Value = 103.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 311.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -81.25
$ws.Range("N31").ClearContents()

$ws.Range("H38").Value = 380.72223
$ws.Range("I38").Value = 65.25
$ws.Range("J38").Value = 2904.5
$ws.Range("K38").Value = 195.75
$ws.Range("L38").Value = 8713.5
$ws.Range("M38").Value = 176.25
$ws.Range("N38").Value = -9457.5

$ws.Range("H39").Value = 2433.5
$ws.Range("I39").Value = 479.2
$ws.Range("J39").Value = 7319.25
$ws.Range("K39").Value = 1437.6
$ws.Range("L39").Value = 21957.75
$ws.Range("M39").Value = -1141.6
$ws.Range("N39").Value = -22549.75

$ws.Range("H70").Value = 334283
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 500724.5
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 1502173.5
$ws.Range("M70").Value = -3930
$ws.Range("N70").Value = -1502713.5

$ws.Range("H73").Value = 334283
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 500724.5
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 1502173.5
$ws.Range("M73").Value = -3264
$ws.Range("N73").Value = -1504045.5

$ws.Range("H100").Value = 1504.6364
$ws.Range("I100").Value = 1450.5
$ws.Range("J100").Value = 1535.5714
$ws.Range("K100").Value = 1450.5
$ws.Range("L100").Value = 1535.5714
$ws.Range("M100").Value = -909.5
$ws.Range("N100").Value = -2617.5714

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H137").Value = 2364.658
$ws.Range("J137").Value = 3018.5
$ws.Range("L137").Value = 9055.5
$ws.Range("N137").Value = -14155.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1081.4706
$ws.Range("I2").Value = 1033.5745
$ws.Range("K2").Value = 1033.5745
$ws.Range("M2").Value = -920.5744999999999

$ws.Range("H74").Value = 3203.75
$ws.Range("I74").Value = 2995.6
$ws.Range("J74").Value = 3352.4285
$ws.Range("K74").Value = 2995.6
$ws.Range("L74").Value = 3352.4285
$ws.Range("M74").Value = -2121.6
$ws.Range("N74").Value = -5100.4285

$ws.Range("H77").Value = 3203.75
$ws.Range("I77").Value = 2995.6
$ws.Range("J77").Value = 3352.4285
$ws.Range("K77").Value = 14978
$ws.Range("L77").Value = 16762.1425
$ws.Range("M77").Value = -10610
$ws.Range("N77").Value = -25498.1425

$ws.Range("H107").Value = 55599.6
$ws.Range("J107").Value = 55599.6
$ws.Range("L107").Value = 55599.6
$ws.Range("N107").Value = -63279.6

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H116").Value = 1081.4706
$ws.Range("I116").Value = 1033.5745
$ws.Range("K116").Value = 1033.5745
$ws.Range("M116").Value = 1260.4255

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1081.4706
$ws.Range("I3").Value = 1033.5745
$ws.Range("K3").Value = 1033.5745
$ws.Range("M3").Value = -919.5744999999999

$ws.Range("H20").Value = 5166.5264
$ws.Range("I20").Value = 5122.75
$ws.Range("K20").Value = 5122.75
$ws.Range("M20").Value = -4875.75

$ws.Range("H99").Value = 4690.8335
$ws.Range("I99").Value = 4690.8335
$ws.Range("K99").Value = 4690.8335
$ws.Range("M99").Value = -3192.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 19124.875
$ws.Range("I4").Value = 18999.857
$ws.Range("K4").Value = 18999.857
$ws.Range("M4").Value = -18887.857

$ws.Range("H20").Value = 83000
$ws.Range("J20").Value = 83000
$ws.Range("L20").Value = 83000
$ws.Range("N20").Value = -83472

$ws.Range("H30").Value = 83000
$ws.Range("J30").Value = 83000
$ws.Range("L30").Value = 83000
$ws.Range("N30").Value = -83182

$ws.Range("H99").Value = 6047.826
$ws.Range("I99").Value = 4886.75
$ws.Range("K99").Value = 4886.75
$ws.Range("M99").Value = -3388.75

$ws.Range("H122").Value = 15249.25
$ws.Range("I122").Value = 14998
$ws.Range("J122").Value = 15333
$ws.Range("K122").Value = 44994
$ws.Range("L122").Value = 45999
$ws.Range("M122").Value = -42544
$ws.Range("N122").Value = -50899

$ws.Range("H126").Value = 6047.826
$ws.Range("I126").Value = 4886.75
$ws.Range("K126").Value = 14660.25
$ws.Range("M126").Value = -12190.25

$ws.Range("H128").Value = 83000
$ws.Range("J128").Value = 83000
$ws.Range("L128").Value = 83000
$ws.Range("N128").Value = -92960

$ws.Range("H129").Value = 97598.2
$ws.Range("J129").Value = 97598.2
$ws.Range("L129").Value = 97598.2
$ws.Range("N129").Value = -107598.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 183
$ws.Range("I6").Value = 210.66667
$ws.Range("K6").Value = 632.00001
$ws.Range("M6").Value = -519.00001

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H137").Value = 3716
$ws.Range("J137").Value = 3087.0715
$ws.Range("L137").Value = 9261.2145
$ws.Range("N137").Value = -19461.2145

$ws.Range("H140").Value = 2006.8077
$ws.Range("I140").Value = 1167.2273
$ws.Range("K140").Value = 3501.6819
$ws.Range("M140").Value = 1678.3181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7058.952
$ws.Range("I40").Value = 7058.952
$ws.Range("K40").Value = 7058.952
$ws.Range("M40").Value = -6922.952

$ws.Range("H122").Value = 3975.6667
$ws.Range("I122").Value = 3097
$ws.Range("K122").Value = 9291
$ws.Range("M122").Value = -6841

$ws.Range("H136").Value = 3010.7334
$ws.Range("I136").Value = 3059.1892
$ws.Range("K136").Value = 9177.567599999998
$ws.Range("M136").Value = -6627.567599999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 186029.42
$ws.Range("J5").Value = 186029.42
$ws.Range("L5").Value = 186029.42
$ws.Range("N5").Value = -186253.42

$ws.Range("H96").Value = 1301.375
$ws.Range("I96").Value = 1082.4
$ws.Range("J96").Value = 1666.3334
$ws.Range("K96").Value = 1082.4
$ws.Range("L96").Value = 1666.3334
$ws.Range("M96").Value = 290.5999999999999
$ws.Range("N96").Value = -4412.3334

$ws.Range("H128").Value = 96583.164
$ws.Range("J128").Value = 96583.164
$ws.Range("L128").Value = 96583.164
$ws.Range("N128").Value = -106543.164
